# Commit: "changed tests to simple fcs file"
#
# The "Samples" sheet listed two flow-cytometry samples (A1 and A2) pointing
# at sample001.fcs / sample006.fcs with channel list "FSC,SSC,FL1". The
# author trimmed this down to a single sample (A1) that points at a small
# test fixture (small.fcs) with an updated, more realistic channel list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# Update row 2 (sample A1): new data location + new channel list.
$ws.Range("E2").Value = "/Users/qr24461/OneDrive - University of Bristol/Code/eebio-tools/test/inputs/small.fcs"
$ws.Range("F2").Value = "FSC-H,SSC-H,FL1-H,FL1-H,FL3-H,FL1-A,FL4-H"

# Remove row 3 (sample A2) entirely - only one sample remains.
$ws.Rows("3:3").Delete()

# Make "Samples" the active sheet/selection, with E3 selected there.
$ws.Activate()
$ws.Range("E3").Select()
